# Mensaje descriptivo de los cambios
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 held "ENF-777-179" -> becomes "ENF-777-555"
$ws.Range("A2").Value = "ENF-777-555"

# A3 held "MED-777-180" -> stays the same text (its shared-string slot just
# gets reused after the old ENF-777-179 entry is replaced), set explicitly
# to keep the written value consistent with the source data.
$ws.Range("A3").Value = "MED-777-180"

# Update the active selection to A7
$ws.Range("A7").Select()
